$d = $word.ActiveDocument

# --- Change 1: remove the "Meta description: ..." paragraph that
# currently follows the H1 title paragraph. ---
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# --- Change 2: insert a new bold paragraph
# "Play Big Fin Bay Slot for Free - Exciting Marlin Hunting Theme"
# right before the final ("Prompt: ...") paragraph. ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insPoint = $lastPara.Range.Start
$insRange = $d.Range($insPoint, $insPoint)

$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Big Fin Bay Slot for Free - Exciting Marlin Hunting Theme</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$insRange.InsertXML($xmlFrag)

# InsertXML leaves a stray empty paragraph behind the one we just
# created (it split the destination paragraph); remove it.
$strayPara = $d.Paragraphs.Item($count + 1)
if ($strayPara.Range.Text.Trim() -eq "") {
    $strayPara.Range.Delete()
}

# --- Change 3: replace the "Prompt: ..." text with the new
# meta-description copy, keeping the italic run formatting intact. ---
$oldText = 'Prompt: Create a feature image for "Big Fin Bay" that showcases the adventure of the old sea wolf and the Marlin, while still reflecting the overall cartoon-style theme of the game. The image should feature the game''s main character, a happy Maya warrior with glasses, as he sets sail on his fishing boat with an optimistic grin on his face. The background should depict the beautiful, enchanting bay, with glimpses of the Marlin jumping out of the water. The image should be bright and colorful, with a playful tone that captures the excitement and thrill of the game.'
$newText = 'Read our review of Big Fin Bay Slot and play for free with an expanding Wild symbol, free spins, and 117,649 ways to win during bonus rounds.'

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
